# BUG instruction images not showing
# - Remove stray empty placeholder cells (C3:D5) left over on RatingPreCond
# - Add a new "PreCond1" worksheet (a second PreCond-style results sheet) with its data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. RatingPreCond: drop the leftover empty inlineStr cells in C3,D3 / C4,D4 / C5,D5
# ---------------------------------------------------------------------------
$wsRating = $wb.Worksheets.Item("RatingPreCond")
$wsRating.Range("C3:D5").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new "PreCond1" worksheet after "RatingPreCond"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "PreCond1"

# --- header row ---
$ws.Range("A1").Value = "PreCondName"
$ws.Range("B1").Value = "TrgCol"
$ws.Range("C1").Value = "n"
$ws.Range("D1").Value = "PreCondTriggKey.keys_raw"
$ws.Range("E1").Value = "PreCondTriggMouse.leftButton_raw"
$ws.Range("F1").Value = "PreCondTriggMouse.midButton_raw"
$ws.Range("G1").Value = "PreCondTriggMouse.rightButton_raw"
$ws.Range("H1").Value = "PreCondTriggMouse.time_raw"
$ws.Range("I1").Value = "PreCondTriggMouse.x_raw"
$ws.Range("J1").Value = "PreCondTriggMouse.y_raw"
$ws.Range("K1").Value = "order"

# --- column A: stimulus name ---
$ws.Range("A2:A6").Value = "stimuli/Neg.BMP"
$ws.Range("A7:A11").Value = "stimuli/Neu.BMP"
$ws.Range("A12:A16").Value = "stimuli/Pos.BMP"
$ws.Range("A17:A19").Value = "stimuli/Trig.BMP"

# --- column B: TrgCol ---
$ws.Range("B2:B6").Value = 0
$ws.Range("B7:B11").Value = 0.4
$ws.Range("B12:B16").Value = 0.6
$ws.Range("B17:B19").Value = 0.9

# --- column C: n ---
$ws.Range("C2:C16").Value = 0
$ws.Range("C17").Value = 1
$ws.Range("C18:C19").Value = 0

# --- column E17: PreCondTriggMouse.leftButton_raw for the one keypress trial ---
$ws.Range("E17").Value = 0

# --- columns D:J, rows 2-16 & 18-19: literal '--' placeholder text ---
# (uses a helper formula so the leading apostrophe is kept as literal text
#  instead of being treated as Excel's "quote prefix" marker)
$ws.Range("D2:J16").Formula = "=CHAR(39)&""--""&CHAR(39)"
$ws.Range("D18:J19").Formula = "=CHAR(39)&""--""&CHAR(39)"
$ws.Range("D2:J19").Copy() | Out-Null
$ws.Range("D2:J19").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# --- extraInfo block ---
$ws.Range("A21").Value = "extraInfo"
$ws.Range("A22").Value = "Participant_ID"
$ws.Range("B22").Value = "S00dsd"
$ws.Range("A23").Value = "Version"
$ws.Range("B23").Value = "Short"
$ws.Range("A24").Value = "Language"
$ws.Range("B24").Value = "EN"
$ws.Range("A25").Value = "date"
$ws.Range("B25").Value = "2023-06-28_00h05.04.770"
$ws.Range("A26").Value = "expName"
$ws.Range("B26").Value = "TCET"
$ws.Range("A27").Value = "psychopyVersion"
$ws.Range("A28").Value = "frameRate"
$ws.Range("B28").Value = 59.96296686889758

# B27 ("2023.1.1") would be auto-parsed as a date by a plain .Value assignment,
# so build it as a text formula result and paste back as a literal value instead.
$ws.Range("B27").Formula = "=""2023.1.1"""
$ws.Range("B27").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
